$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("20:20").Insert()
$ws.Cells.Item(20,1).EntireRow.RowHeight = 120
